$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new row of data (this appends shared string "地面容器" and values)
$ws.Range("A7").Value = 110000
$ws.Range("B7").Value = "地面容器"
$ws.Range("C7").Value = 99999999

# Normalize cell styles: unify all data/header cells onto a single style
# (matches the cellXfs collapsing from 4 entries down to 2 in the diff).
$dataRange = $ws.Range("A1:C7")
$dataRange.Style = "Normal"
$dataRange.Font.Name = "HarmonyOS Sans SC"
$headerRange = $ws.Range("A1:C1")
$headerRange.HorizontalAlignment = -4108
$dataRange.VerticalAlignment = -4108

# Column width: column A narrower than the rest
$ws.Columns.Item(1).ColumnWidth = 9.6

# Update selection to match target state
$ws.Range("D9").Select()
